$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B column (sections) for rows 4-13
$ws.Range("B4").Value = "4_CSE_B, 4_CSE_C, 4_CSE_E, 4_AIDS_H, 4_CSBT_I, 4_IT_K, 4_ICT_M, 4_CSBS_N, 6_CSE_B, 6_CSE_D, 6_CSE_E, 6_AIDS_F, 6_AIDS_H, 6_IT_L"
$ws.Range("B5").Value = "4_CSE_A, 4_CSE_D, 4_AIDS_F, 4_AIDS_G, 4_IoTA_J, 4_IT_L, 6_CSE_A, 6_CSE_C, 6_AIDS_G, 6_CSBT_I, 6_IoTA_J, 6_IT_K, 6_ICT_M, 6_CSBS_N"
$ws.Range("B6").Value = "4_CSE_B, 4_AIDS_F, 4_AIDS_G, 4_AIDS_H, 4_IT_K, 4_CSBS_N, 6_CSE_B, 6_CSE_D, 6_CSE_E, 6_AIDS_F, 6_AIDS_G, 6_CSBT_I, 6_IoTA_J, 6_IT_K"
$ws.Range("B7").Value = "4_CSE_A, 4_CSE_C, 4_CSE_D, 4_CSE_E, 4_CSBT_I, 4_IoTA_J, 4_IT_L, 4_ICT_M, 6_CSE_A, 6_CSE_C, 6_AIDS_H, 6_IT_L, 6_ICT_M, 6_CSBS_N"
$ws.Range("B8").Value = "4_CSE_B, 4_CSE_C, 4_CSE_E, 4_AIDS_F, 4_AIDS_G, 4_CSBT_I, 4_IT_K, 4_IT_L, 4_CSBS_N, 6_CSE_E, 6_IoTA_J, 6_IT_K, 6_IT_L, 6_CSBS_N"
$ws.Range("B9").Value = "4_CSE_A, 4_CSE_D, 4_AIDS_H, 4_IoTA_J, 4_ICT_M, 6_CSE_A, 6_CSE_B, 6_CSE_C, 6_CSE_D, 6_AIDS_F, 6_AIDS_G, 6_AIDS_H, 6_CSBT_I, 6_ICT_M"
$ws.Range("B10").Value = "4_CSE_A, 4_CSE_B, 4_CSE_C, 4_CSE_D, 4_AIDS_F, 4_IoTA_J, 4_IT_L, 4_ICT_M, 6_CSE_A, 6_CSE_C, 6_CSBT_I, 6_IoTA_J, 6_ICT_M, 6_CSBS_N"
$ws.Range("B11").Value = "4_CSE_E, 4_AIDS_G, 4_AIDS_H, 4_CSBT_I, 4_IT_K, 4_CSBS_N, 6_CSE_B, 6_CSE_D, 6_CSE_E, 6_AIDS_F, 6_AIDS_G, 6_AIDS_H, 6_IT_K, 6_IT_L"
$ws.Range("B12").Value = "4_CSE_A, 4_CSE_D, 4_AIDS_G, 4_AIDS_H, 4_CSBT_I, 4_IoTA_J, 4_IT_L, 6_CSE_A, 6_CSE_C, 6_AIDS_F, 6_AIDS_G, 6_AIDS_H, 6_CSBT_I, 6_ICT_M"
$ws.Range("B13").Value = "4_CSE_B, 4_CSE_C, 4_CSE_E, 4_AIDS_F, 4_IT_K, 4_ICT_M, 4_CSBS_N, 6_CSE_B, 6_CSE_D, 6_CSE_E, 6_IoTA_J, 6_IT_K, 6_IT_L, 6_CSBS_N"

# Update G column (block) for rows 14-99
$ws.Range("G14").Value = "(2, 6), (2, 7)"
$ws.Range("G15").Value = "(2, 0), (2, 1)"
$ws.Range("G16").Value = "(1, 6), (1, 7)"
$ws.Range("G17").Value = "(3, 6), (3, 7)"
$ws.Range("G18").Value = "(0, 0), (0, 1)"
$ws.Range("G19").Value = "(1, 0), (1, 1)"
$ws.Range("G20").Value = "(4, 2), (4, 3)"
$ws.Range("G21").Value = "(3, 2), (3, 3)"
$ws.Range("G22").Value = "(0, 2), (0, 3)"
$ws.Range("G23").Value = "(1, 6), (1, 7)"
$ws.Range("G24").Value = "(4, 6), (4, 7)"
$ws.Range("G25").Value = "(3, 6), (3, 7)"
$ws.Range("G27").Value = "(4, 6), (4, 7)"
$ws.Range("G28").Value = "(3, 0), (3, 1)"
$ws.Range("G29").Value = "(3, 2), (3, 3)"
$ws.Range("G30").Value = "(4, 6), (4, 7)"
$ws.Range("G31").Value = "(0, 4), (0, 5)"
$ws.Range("G32").Value = "(1, 6), (1, 7)"
$ws.Range("G33").Value = "(4, 6), (4, 7)"
$ws.Range("G34").Value = "(1, 4), (1, 5)"
$ws.Range("G35").Value = "(3, 6), (3, 7)"
$ws.Range("G36").Value = "(0, 0), (0, 1)"
$ws.Range("G37").Value = "(1, 0), (1, 1)"
$ws.Range("G38").Value = "(3, 6), (3, 7)"
$ws.Range("G39").Value = "(4, 6), (4, 7)"
$ws.Range("G41").Value = "(2, 6), (2, 7)"
$ws.Range("G42").Value = "(0, 6), (0, 7)"
$ws.Range("G43").Value = "(3, 6), (3, 7)"
$ws.Range("G44").Value = "(3, 6), (3, 7)"
$ws.Range("G45").Value = "(2, 6), (2, 7)"
$ws.Range("G46").Value = "(0, 6), (0, 7)"
$ws.Range("G47").Value = "(0, 0), (0, 1)"
$ws.Range("G48").Value = "(4, 4), (4, 5)"
$ws.Range("G49").Value = "(3, 6), (3, 7)"
$ws.Range("G50").Value = "(2, 6), (2, 7)"
$ws.Range("G51").Value = "(4, 4), (4, 5)"
$ws.Range("G52").Value = "(3, 6), (3, 7)"
$ws.Range("G53").Value = "(4, 6), (4, 7)"
$ws.Range("G54").Value = "(2, 2), (2, 3)"
$ws.Range("G55").Value = "(3, 6), (3, 7)"
$ws.Range("G56").Value = "(0, 2), (0, 3)"
$ws.Range("G57").Value = "(3, 6), (3, 7)"
$ws.Range("G58").Value = "(4, 4), (4, 5)"
$ws.Range("G59").Value = "(4, 6), (4, 7)"
$ws.Range("G60").Value = "(2, 2), (2, 3)"
$ws.Range("G61").Value = "(0, 0), (0, 1)"
$ws.Range("G62").Value = "(2, 4), (2, 5)"
$ws.Range("G63").Value = "(0, 0), (0, 1)"
$ws.Range("G64").Value = "(3, 2), (3, 3)"
$ws.Range("G65").Value = "(3, 0), (3, 1)"
$ws.Range("G66").Value = "(2, 6), (2, 7)"
$ws.Range("G67").Value = "(0, 6), (0, 7)"
$ws.Range("G68").Value = "(3, 0), (3, 1)"
$ws.Range("G69").Value = "(2, 6), (2, 7)"
$ws.Range("G70").Value = "(1, 2), (1, 3)"
$ws.Range("G71").Value = "(2, 2), (2, 3)"
$ws.Range("G72").Value = "(3, 4), (3, 5)"
$ws.Range("G73").Value = "(0, 4), (0, 5)"
$ws.Range("G74").Value = "(3, 2), (3, 3)"
$ws.Range("G75").Value = "(1, 6), (1, 7)"
$ws.Range("G76").Value = "(2, 6), (2, 7)"
$ws.Range("G77").Value = "(0, 0), (0, 1)"
$ws.Range("G78").Value = "(4, 0), (4, 1)"
$ws.Range("G79").Value = "(0, 4), (0, 5)"
$ws.Range("G80").Value = "(3, 2), (3, 3)"
$ws.Range("G81").Value = "(1, 0), (1, 1)"
$ws.Range("G82").Value = "(0, 0), (0, 1)"
$ws.Range("G83").Value = "(4, 4), (4, 5)"
$ws.Range("G84").Value = "(3, 4), (3, 5)"
$ws.Range("G85").Value = "(0, 0), (0, 1)"
$ws.Range("G86").Value = "(2, 4), (2, 5)"
$ws.Range("G87").Value = "(4, 0), (4, 1)"
$ws.Range("G88").Value = "(3, 4), (3, 5)"
$ws.Range("G89").Value = "(1, 0), (1, 1)"
$ws.Range("G90").Value = "(2, 6), (2, 7)"
$ws.Range("G91").Value = "(3, 2), (3, 3)"
$ws.Range("G92").Value = "(4, 2), (4, 3)"
$ws.Range("G93").Value = "(0, 4), (0, 5)"
$ws.Range("G94").Value = "(0, 0), (0, 1)"
$ws.Range("G95").Value = "(4, 2), (4, 3)"
$ws.Range("G96").Value = "(2, 4), (2, 5)"
$ws.Range("G97").Value = "(1, 0), (1, 1)"
$ws.Range("G98").Value = "(1, 0), (1, 1)"
$ws.Range("G99").Value = "(0, 2), (0, 3)"
